$d = $word.ActiveDocument

# Replace the placeholder id text (and eat the trailing space run) in the
# first paragraph in one shot.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5301_topic_28__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5301_601_91__ID**", 2)

# Give that paragraph a (zero-weight) paragraph border with 5pt spacing on
# every side, and widen its left indent from 120 to 225 twips.
$p = $d.Paragraphs(1)
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
$p.Format.LeftIndent = 11.25
